$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update November 2025 row (row 24) stats
$ws.Range("B24").Value = 6386
$ws.Range("C24").Value = 1003
$ws.Range("D24").Value = 5965833
$ws.Range("E24").Value = 934.204979642969
$ws.Range("F24").Value = 8.864643709512453
$ws.Range("G24").Value = 3.937823834196896
$ws.Range("H24").Value = 26.37945147325367
